$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to be stored as text so exact formatting
# (trailing zeros, thousand-dot separators, etc.) survives the write,
# matching the source data which stores prices as inline strings.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.099.60"
$ws.Range("E2").Value = "  +5.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.782.38"
$ws.Range("E3").Value = "  +3.23%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.17"
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9991"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4915"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2679"
$ws.Range("E8").Value = "  +2.20%  "
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.779.29"
$ws.Range("E10").Value = "  +2.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.41"
$ws.Range("E11").Value = "  +3.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07040"
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6283"
$ws.Range("E13").Value = "  +2.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.654"
$ws.Range("E14").Value = "  +3.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "80.23"
$ws.Range("E15").Value = "  +3.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.090.19"
$ws.Range("E16").Value = "  +5.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9981"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9974"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007231"
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("E20").Value = "  +4.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.011.42"
$ws.Range("E21").Value = "  +3.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.578"
$ws.Range("E22").Value = "  +2.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.715"
$ws.Range("E23").Value = "  +1.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.233"
$ws.Range("E24").Value = "  +2.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.58"
$ws.Range("E25").Value = "  +2.53%  "
$ws.Range("E26").Value = "  +2.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.862"
$ws.Range("E27").Value = "  +5.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "109.14"
$ws.Range("E28").Value = "  +2.45%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.181"
$ws.Range("E30").Value = "  +6.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08281"
$ws.Range("E31").Value = "  +3.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.794"
$ws.Range("E32").Value = "  +3.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04893"
$ws.Range("E33").Value = "  +8.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.075"
$ws.Range("E34").Value = "  +7.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6542"
$ws.Range("E35").Value = "  +4.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.608"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9537"
$ws.Range("E37").Value = "  +2.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.585"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.056"
$ws.Range("E39").Value = "  +1.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.973"
$ws.Range("E40").Value = "  +7.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01553"
$ws.Range("E41").Value = "  +2.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9989"
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.90"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3987"
$ws.Range("E44").Value = "  +3.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.178"
$ws.Range("E45").Value = "  +3.68%  "
$ws.Range("E46").Value = "  +3.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05428"
$ws.Range("E47").Value = "  +0.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.002"
$ws.Range("E48").Value = "  +1.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.300"
$ws.Range("E49").Value = "  +5.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.64"
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.92"
$ws.Range("E51").Value = "  +2.27%  "
